$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value2 = 6614.5
$ws.Range("I62").Value2 = 4004.6667
$ws.Range("J62").Value2 = 14444
$ws.Range("K62").Value2 = 4004.6667
$ws.Range("L62").Value2 = 14444
$ws.Range("M62").Value2 = -3380.6667
$ws.Range("N62").Value2 = -15692
$ws.Range("H65").Value2 = 6614.5
$ws.Range("I65").Value2 = 4004.6667
$ws.Range("J65").Value2 = 14444
$ws.Range("K65").Value2 = 20023.3335
$ws.Range("L65").Value2 = 72220
$ws.Range("M65").Value2 = -16903.3335
$ws.Range("N65").Value2 = -78460
$ws.Range("H70").Value2 = 3100
$ws.Range("J70").Value2 = 3200
$ws.Range("L70").Value2 = 9600
$ws.Range("N70").Value2 = -10140
$ws.Range("H73").Value2 = 3100
$ws.Range("J73").Value2 = 3200
$ws.Range("L73").Value2 = 9600
$ws.Range("N73").Value2 = -11472

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value2 = 3973.3
$ws.Range("I45").Value2 = 2503.8
$ws.Range("K45").Value2 = 2503.8
$ws.Range("M45").Value2 = -2126.8
$ws.Range("H88").Value2 = 2601.75
$ws.Range("I88").Value2 = 2194.111
$ws.Range("J88").Value2 = 2935.2727
$ws.Range("K88").Value2 = 2194.111
$ws.Range("L88").Value2 = 2935.2727
$ws.Range("M88").Value2 = -1788.111
$ws.Range("N88").Value2 = -3747.2727
$ws.Range("H91").Value2 = 2601.75
$ws.Range("I91").Value2 = 2194.111
$ws.Range("J91").Value2 = 2935.2727
$ws.Range("K91").Value2 = 2194.111
$ws.Range("L91").Value2 = 2935.2727
$ws.Range("M91").Value2 = -790.1109999999999
$ws.Range("N91").Value2 = -5743.2727
$ws.Range("H97").Value2 = 99.2
$ws.Range("I97").Value2 = 126.25
$ws.Range("J97").Value2 = 81.166664
$ws.Range("K97").Value2 = 126.25
$ws.Range("L97").Value2 = 81.166664
$ws.Range("M97").Value2 = 369.75
$ws.Range("N97").Value2 = -1073.166664
$ws.Range("H132").Value2 = 1230.9231
$ws.Range("I132").Value2 = 1148.3334
$ws.Range("J132").Value2 = 2222
$ws.Range("K132").Value2 = 3445.0002
$ws.Range("L132").Value2 = 6666
$ws.Range("M132").Value2 = -915.0001999999999
$ws.Range("N132").Value2 = -11726

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value2 = 7197.7144
$ws.Range("I54").Value2 = 5398.1665
$ws.Range("J54").Value2 = 17995
$ws.Range("K54").Value2 = 5398.1665
$ws.Range("L54").Value2 = 17995
$ws.Range("M54").Value2 = -4914.1665
$ws.Range("N54").Value2 = -18963
$ws.Range("H99").Value2 = 2574
$ws.Range("I99").Value2 = 2430.4546
$ws.Range("K99").Value2 = 2430.4546
$ws.Range("M99").Value2 = -932.4546
$ws.Range("H134").Value2 = 5999.75
$ws.Range("I134").Value2 = 5999
$ws.Range("K134").Value2 = 17997
$ws.Range("M134").Value2 = -15462

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 1512.7142
$ws.Range("I22").Value2 = 944.5
$ws.Range("J22").Value2 = 1740
$ws.Range("K22").Value2 = 944.5
$ws.Range("L22").Value2 = 1740
$ws.Range("M22").Value2 = -594.5
$ws.Range("N22").Value2 = -2440
$ws.Range("H25").Value2 = 570
$ws.Range("I25").Value2 = 570
$ws.Range("K25").Value2 = 570
$ws.Range("M25").Value2 = -396
$ws.Range("H94").Value2 = 2801
$ws.Range("I94").Value2 = 3010.5454
$ws.Range("K94").Value2 = 3010.5454
$ws.Range("M94").Value2 = -2559.5454
$ws.Range("H107").Value2 = 502.6
$ws.Range("I107").Value2 = 680
$ws.Range("J107").Value2 = 384.33334
$ws.Range("K107").Value2 = 680
$ws.Range("L107").Value2 = 384.33334
$ws.Range("M107").Value2 = 1240
$ws.Range("N107").Value2 = -4224.33334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value2 = 255.2
$ws.Range("I50").Value2 = 255.2
$ws.Range("J50").Value2 = 0
$ws.Range("K50").Value2 = 765.5999999999999
$ws.Range("L50").Value2 = 0
$ws.Range("M50").Value2 = -284.5999999999999
$ws.Range("N50").ClearContents()
$ws.Range("H53").Value2 = 255.2
$ws.Range("I53").Value2 = 255.2
$ws.Range("J53").Value2 = 0
$ws.Range("K53").Value2 = 765.5999999999999
$ws.Range("L53").Value2 = 0
$ws.Range("M53").Value2 = -284.5999999999999
$ws.Range("N53").ClearContents()
$ws.Range("H94").Value2 = 1000
$ws.Range("I94").Value2 = 1000
$ws.Range("K94").Value2 = 3000
$ws.Range("M94").Value2 = -2324
$ws.Range("H109").Value2 = 332
$ws.Range("J109").Value2 = 0
$ws.Range("L109").Value2 = 0
$ws.Range("N109").ClearContents()
$ws.Range("H136").Value2 = 12998
$ws.Range("I136").Value2 = 12998
$ws.Range("K136").Value2 = 38994
$ws.Range("M136").Value2 = -33894

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value2 = 433.48148
$ws.Range("I97").Value2 = 408.14285
$ws.Range("J97").Value2 = 522.1667
$ws.Range("K97").Value2 = 408.14285
$ws.Range("L97").Value2 = 522.1667
$ws.Range("M97").Value2 = 87.85714999999999
$ws.Range("N97").Value2 = -1514.1667
$ws.Range("H139").Value2 = 90000
$ws.Range("J139").Value2 = 90000
$ws.Range("L139").Value2 = 90000
$ws.Range("N139").Value2 = -100280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 952.1667
$ws.Range("J22").Value2 = 1081.6666
$ws.Range("L22").Value2 = 1081.6666
$ws.Range("N22").Value2 = -1671.6666
$ws.Range("H27").Value2 = 952.1667
$ws.Range("J27").Value2 = 1081.6666
$ws.Range("L27").Value2 = 1081.6666
$ws.Range("N27").Value2 = -1295.6666
$ws.Range("H40").Value2 = 1707.8572
$ws.Range("I40").Value2 = 1707.8572
$ws.Range("J40").Value2 = 0
$ws.Range("K40").Value2 = 1707.8572
$ws.Range("L40").Value2 = 0
$ws.Range("M40").Value2 = -1571.8572
$ws.Range("N40").ClearContents()
$ws.Range("H93").Value2 = 0
$ws.Range("I93").Value2 = 0
$ws.Range("K93").Value2 = 0
$ws.Range("M93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value2 = 0
$ws.Range("I40").Value2 = 0
$ws.Range("K40").Value2 = 0
$ws.Range("M40").ClearContents()
$ws.Range("H107").Value2 = 410.3
$ws.Range("I107").Value2 = 441.57144
$ws.Range("J107").Value2 = 337.33334
$ws.Range("K107").Value2 = 1324.71432
$ws.Range("L107").Value2 = 1012.00002
$ws.Range("M107").Value2 = 595.28568
$ws.Range("N107").Value2 = -4852.00002
$ws.Range("H113").Value2 = 900
$ws.Range("I113").Value2 = 900
$ws.Range("K113").Value2 = 2700
$ws.Range("M113").Value2 = -530
$ws.Range("H122").Value2 = 492.7
$ws.Range("I122").Value2 = 500.22223
$ws.Range("K122").Value2 = 1500.66669
$ws.Range("M122").Value2 = 949.33331
